$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 112.8140721321106
$ws.Range("C2").Value = 1.151494228068786
$ws.Range("D2").Value = 1.034724760055542
$ws.Range("E2").Value = 0.02641201518394986
$ws.Range("J2").Value = 0.7489361702127659
$ws.Range("K2").Value = 0.7021276595744681
$ws.Range("L2").Value = 0.7446808510638298
$ws.Range("M2").Value = 0.7446808510638298
$ws.Range("N2").Value = 0.7521367521367521
$ws.Range("O2").Value = 0.7385124568103292
$ws.Range("P2").Value = 0.0184077530919213
$ws.Range("B3").Value = 222.0525091648102
$ws.Range("C3").Value = 1.048401361047337
$ws.Range("D3").Value = 1.07569899559021
$ws.Range("E3").Value = 0.05819221593704209
$ws.Range("J3").Value = 0.7829787234042553
$ws.Range("K3").Value = 0.7617021276595745
$ws.Range("L3").Value = 0.7489361702127659
$ws.Range("M3").Value = 0.7574468085106383
$ws.Range("N3").Value = 0.782051282051282
$ws.Range("O3").Value = 0.7666230223677032
$ws.Range("P3").Value = 0.01361457305897402
$ws.Range("B4").Value = 440.5179897785187
$ws.Range("C4").Value = 1.958768132014467
$ws.Range("D4").Value = 1.081076097488403
$ws.Range("E4").Value = 0.0515729379102302
$ws.Range("J4").Value = 0.8
$ws.Range("K4").Value = 0.8
$ws.Range("L4").Value = 0.8085106382978723
$ws.Range("M4").Value = 0.8127659574468085
$ws.Range("N4").Value = 0.8418803418803419
$ws.Range("O4").Value = 0.8126313875250044
$ws.Range("P4").Value = 0.01543764176129924
$ws.Range("Q4").Value = 7
$ws.Range("B5").Value = 113.4469215393066
$ws.Range("C5").Value = 1.167578068962321
$ws.Range("D5").Value = 1.108101654052734
$ws.Range("E5").Value = 0.111818747880678
$ws.Range("J5").Value = 0.7489361702127659
$ws.Range("K5").Value = 0.7021276595744681
$ws.Range("L5").Value = 0.7446808510638298
$ws.Range("M5").Value = 0.7446808510638298
$ws.Range("N5").Value = 0.7521367521367521
$ws.Range("O5").Value = 0.7385124568103292
$ws.Range("P5").Value = 0.0184077530919213
$ws.Range("B6").Value = 222.7793167591095
$ws.Range("C6").Value = 1.667646333996571
$ws.Range("D6").Value = 1.202676725387573
$ws.Range("E6").Value = 0.1552925794604522
$ws.Range("J6").Value = 0.7829787234042553
$ws.Range("K6").Value = 0.7617021276595745
$ws.Range("L6").Value = 0.7489361702127659
$ws.Range("M6").Value = 0.7574468085106383
$ws.Range("N6").Value = 0.782051282051282
$ws.Range("O6").Value = 0.7666230223677032
$ws.Range("P6").Value = 0.01361457305897402
$ws.Range("B7").Value = 444.4782149791718
$ws.Range("C7").Value = 1.370579460449444
$ws.Range("D7").Value = 1.176820468902588
$ws.Range("E7").Value = 0.1242808122385996
$ws.Range("J7").Value = 0.8
$ws.Range("K7").Value = 0.8
$ws.Range("L7").Value = 0.8085106382978723
$ws.Range("M7").Value = 0.8127659574468085
$ws.Range("N7").Value = 0.8418803418803419
$ws.Range("O7").Value = 0.8126313875250044
$ws.Range("P7").Value = 0.01543764176129924
$ws.Range("Q7").Value = 7
$ws.Range("B8").Value = 115.6096325874329
$ws.Range("C8").Value = 0.9964905168043321
$ws.Range("D8").Value = 1.091544628143311
$ws.Range("E8").Value = 0.07620130726905652
$ws.Range("J8").Value = 0.7489361702127659
$ws.Range("K8").Value = 0.7021276595744681
$ws.Range("L8").Value = 0.7446808510638298
$ws.Range("M8").Value = 0.7446808510638298
$ws.Range("N8").Value = 0.7521367521367521
$ws.Range("O8").Value = 0.7385124568103292
$ws.Range("P8").Value = 0.0184077530919213
$ws.Range("B9").Value = 224.1770779132843
$ws.Range("C9").Value = 0.9579008126353754
$ws.Range("D9").Value = 1.134570121765137
$ws.Range("E9").Value = 0.06490457096109295
$ws.Range("J9").Value = 0.7829787234042553
$ws.Range("K9").Value = 0.7617021276595745
$ws.Range("L9").Value = 0.7489361702127659
$ws.Range("M9").Value = 0.7574468085106383
$ws.Range("N9").Value = 0.782051282051282
$ws.Range("O9").Value = 0.7666230223677032
$ws.Range("P9").Value = 0.01361457305897402
$ws.Range("B10").Value = 442.4488591194153
$ws.Range("C10").Value = 2.209128231185954
$ws.Range("D10").Value = 1.133418369293213
$ws.Range("E10").Value = 0.1197213124887748
$ws.Range("J10").Value = 0.8
$ws.Range("K10").Value = 0.8
$ws.Range("L10").Value = 0.8085106382978723
$ws.Range("M10").Value = 0.8127659574468085
$ws.Range("N10").Value = 0.8418803418803419
$ws.Range("O10").Value = 0.8126313875250044
$ws.Range("P10").Value = 0.01543764176129924
$ws.Range("Q10").Value = 7
$ws.Range("B11").Value = 178.8798421859741
$ws.Range("C11").Value = 0.4307670391946885
$ws.Range("D11").Value = 1.134435510635376
$ws.Range("E11").Value = 0.1216022443007243
$ws.Range("J11").Value = 0.8042553191489362
$ws.Range("K11").Value = 0.7872340425531915
$ws.Range("L11").Value = 0.8042553191489362
$ws.Range("M11").Value = 0.8212765957446808
$ws.Range("N11").Value = 0.8333333333333334
$ws.Range("O11").Value = 0.8100709219858157
$ws.Range("P11").Value = 0.01584848535995849
$ws.Range("Q11").Value = 10
$ws.Range("B12").Value = 345.8915525913238
$ws.Range("C12").Value = 0.7469578292970765
$ws.Range("D12").Value = 1.080759525299072
$ws.Range("E12").Value = 0.07770711176241575
$ws.Range("J12").Value = 0.8085106382978723
$ws.Range("K12").Value = 0.8042553191489362
$ws.Range("L12").Value = 0.8212765957446808
$ws.Range("M12").Value = 0.8085106382978723
$ws.Range("N12").Value = 0.8589743589743589
$ws.Range("O12").Value = 0.8203055100927441
$ws.Range("P12").Value = 0.02015970987797284
$ws.Range("B13").Value = 674.1369040966034
$ws.Range("C13").Value = 1.613823345030307
$ws.Range("D13").Value = 1.177210712432861
$ws.Range("E13").Value = 0.221941187905537
$ws.Range("J13").Value = 0.7957446808510639
$ws.Range("K13").Value = 0.8042553191489362
$ws.Range("L13").Value = 0.825531914893617
$ws.Range("M13").Value = 0.8340425531914893
$ws.Range("N13").Value = 0.8547008547008547
$ws.Range("O13").Value = 0.8228550645571922
$ws.Range("P13").Value = 0.0211064123596872
$ws.Range("B14").Value = 176.3877244472504
$ws.Range("C14").Value = 0.5784019296200071
$ws.Range("D14").Value = 1.152410888671875
$ws.Range("E14").Value = 0.1182071519540338
$ws.Range("J14").Value = 0.8042553191489362
$ws.Range("K14").Value = 0.7872340425531915
$ws.Range("L14").Value = 0.8042553191489362
$ws.Range("M14").Value = 0.8212765957446808
$ws.Range("N14").Value = 0.8333333333333334
$ws.Range("O14").Value = 0.8100709219858157
$ws.Range("P14").Value = 0.01584848535995849
$ws.Range("Q14").Value = 10
$ws.Range("B15").Value = 346.1172390937805
$ws.Range("C15").Value = 0.8531804618472837
$ws.Range("D15").Value = 1.12978663444519
$ws.Range("E15").Value = 0.1077581435039006
$ws.Range("J15").Value = 0.8085106382978723
$ws.Range("K15").Value = 0.8042553191489362
$ws.Range("L15").Value = 0.8212765957446808
$ws.Range("M15").Value = 0.8085106382978723
$ws.Range("N15").Value = 0.8589743589743589
$ws.Range("O15").Value = 0.8203055100927441
$ws.Range("P15").Value = 0.02015970987797284
$ws.Range("B16").Value = 673.9129508972168
$ws.Range("C16").Value = 1.242048769936551
$ws.Range("D16").Value = 1.129365062713623
$ws.Range("E16").Value = 0.1592003479618045
$ws.Range("J16").Value = 0.7957446808510639
$ws.Range("K16").Value = 0.8042553191489362
$ws.Range("L16").Value = 0.825531914893617
$ws.Range("M16").Value = 0.8340425531914893
$ws.Range("N16").Value = 0.8547008547008547
$ws.Range("O16").Value = 0.8228550645571922
$ws.Range("P16").Value = 0.0211064123596872
$ws.Range("B17").Value = 176.5527980327606
$ws.Range("C17").Value = 0.9598749726479374
$ws.Range("D17").Value = 1.24395604133606
$ws.Range("E17").Value = 0.1544675446691934
$ws.Range("J17").Value = 0.8042553191489362
$ws.Range("K17").Value = 0.7872340425531915
$ws.Range("L17").Value = 0.8042553191489362
$ws.Range("M17").Value = 0.8212765957446808
$ws.Range("N17").Value = 0.8333333333333334
$ws.Range("O17").Value = 0.8100709219858157
$ws.Range("P17").Value = 0.01584848535995849
$ws.Range("Q17").Value = 10
$ws.Range("B18").Value = 340.591864490509
$ws.Range("C18").Value = 6.641816336182635
$ws.Range("D18").Value = 0.9923890590667724
$ws.Range("E18").Value = 0.0617189675840896
$ws.Range("J18").Value = 0.8085106382978723
$ws.Range("K18").Value = 0.8042553191489362
$ws.Range("L18").Value = 0.8212765957446808
$ws.Range("M18").Value = 0.8085106382978723
$ws.Range("N18").Value = 0.8589743589743589
$ws.Range("O18").Value = 0.8203055100927441
$ws.Range("P18").Value = 0.02015970987797284
$ws.Range("B19").Value = 541.1166676044464
$ws.Range("C19").Value = 22.23276382955041
$ws.Range("D19").Value = 0.7046218395233155
$ws.Range("E19").Value = 0.168336572840162
$ws.Range("J19").Value = 0.7957446808510639
$ws.Range("K19").Value = 0.8042553191489362
$ws.Range("L19").Value = 0.825531914893617
$ws.Range("M19").Value = 0.8340425531914893
$ws.Range("N19").Value = 0.8547008547008547
$ws.Range("O19").Value = 0.8228550645571922
$ws.Range("P19").Value = 0.0211064123596872
